$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.762.31"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").Value = "3.235.79"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.18"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.47"
$ws.Range("E6").Value = "  -3.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("D9").Value = "3.231.69"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.390"
$ws.Range("E12").Value = "  -2.98%  "
$ws.Range("D13").Value = "3.795.01"
$ws.Range("E14").Value = "  -3.06%  "
$ws.Range("D15").Value = "64.877.62"
$ws.Range("E15").Value = "  -1.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.79"
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("D17").Value = "3.228.60"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("E18").Value = "  -2.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "417.23"
$ws.Range("E19").Value = "  -3.96%  "
$ws.Range("E20").Value = "  -2.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.82"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.22"
$ws.Range("E22").Value = "  -2.19%  "
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.51"
$ws.Range("E24").Value = "  -1.58%  "
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("E26").Value = "  +4.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.496"
$ws.Range("E27").Value = "  -1.74%  "
$ws.Range("E28").Value = "  -1.36%  "
$ws.Range("E29").Value = "  +2.05%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("E31").Value = "  -4.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.86"
$ws.Range("E32").Value = "  -1.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.99"
$ws.Range("E34").Value = "  -3.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.43"
$ws.Range("E35").Value = "  -2.67%  "
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.90"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  -2.29%  "
$ws.Range("D39").Value = "2.825.25"
$ws.Range("E39").Value = "  +1.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.44"
$ws.Range("E41").Value = "  -4.36%  "
$ws.Range("E42").Value = "  -2.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.51"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.724"
$ws.Range("E44").Value = "  -6.39%  "
$ws.Range("E45").Value = "  -4.49%  "
$ws.Range("E46").Value = "  -4.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.18"
$ws.Range("E47").Value = "  -4.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "302.51"
$ws.Range("E48").Value = "  -5.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.99"
$ws.Range("E49").Value = "  -5.21%  "
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("E51").Value = "  -1.37%  "
